# Apply the parameter updates described in the commit:
# "Testing different values for analyses to ensure treatment effect
#  modification and standardization to the correct study population."
#
# All of the downstream numbers (potential_preg_trt, the *_preec_* and
# postpreec_preg sheets, etc.) are driven by formulas that reference these
# inputs, so only the raw input cells on potential_preg_untrt need to be
# edited -- everything else recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("potential_preg_untrt")

$ws.Range("C9").Value  = 0.05
$ws.Range("C10").Value = 0.02
$ws.Range("C11").Value = 0.02
$ws.Range("C13").Value = 0.005
$ws.Range("C14").Value = 0.004
$ws.Range("C15").Value = 0.004
$ws.Range("C16").Value = 0.004
$ws.Range("C17").Value = 0.004

# Move the active tab / selection: SimParameters was the selected sheet
# (with C18 selected on potential_preg_untrt); after the edit,
# potential_preg_untrt becomes the active sheet with C12 selected.
$ws.Activate()
$ws.Range("C12").Select()
